$wb = $excel.ActiveWorkbook

# --- ALC sheet updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 7500
$ws.Range("I21").Value = 5000
$ws.Range("K21").Value = 5000
$ws.Range("M21").Value = -4532
$ws.Range("H23").Value = 7500
$ws.Range("I23").Value = 5000
$ws.Range("K23").Value = 5000
$ws.Range("M23").Value = -4766
$ws.Range("H38").Value = 965
$ws.Range("I38").Value = 109.181816
$ws.Range("J38").Value = 1689.1538
$ws.Range("K38").Value = 327.545448
$ws.Range("L38").Value = 5067.4614
$ws.Range("M38").Value = 44.45455200000004
$ws.Range("N38").Value = -5811.4614
$ws.Range("H62").Value = 25013200
$ws.Range("J62").Value = 4900
$ws.Range("L62").Value = 4900
$ws.Range("N62").Value = -6148
$ws.Range("H65").Value = 25013200
$ws.Range("J65").Value = 4900
$ws.Range("L65").Value = 24500
$ws.Range("N65").Value = -30740
$ws.Range("H117").Value = 49743
$ws.Range("J117").Value = 49743
$ws.Range("L117").Value = 49743
$ws.Range("N117").Value = -58921
$ws.Range("H137").Value = 1543.3928
$ws.Range("I137").Value = 1056.4375
$ws.Range("J137").Value = 2192.6667
$ws.Range("K137").Value = 3169.3125
$ws.Range("L137").Value = 6578.000100000001
$ws.Range("M137").Value = -619.3125
$ws.Range("N137").Value = -11678.0001
$ws.Range("H138").Value = 2971.012
$ws.Range("I138").Value = 1462.3823
$ws.Range("J138").Value = 3996.88
$ws.Range("K138").Value = 4387.1469
$ws.Range("L138").Value = 11990.64
$ws.Range("M138").Value = 752.8531000000003
$ws.Range("N138").Value = -22270.64
$ws.Range("H141").Value = 2861.7104
$ws.Range("I141").Value = 1874.697
$ws.Range("J141").Value = 9376
$ws.Range("K141").Value = 5624.090999999999
$ws.Range("L141").Value = 28128
$ws.Range("M141").Value = -444.0909999999994
$ws.Range("N141").Value = -38488

# --- ARM sheet updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 882.1905
$ws.Range("I2").Value = 794.125
$ws.Range("J2").Value = 936.38464
$ws.Range("K2").Value = 794.125
$ws.Range("L2").Value = 936.38464
$ws.Range("M2").Value = -681.125
$ws.Range("N2").Value = -1162.38464
$ws.Range("H110").Value = 2096.1538
$ws.Range("I110").Value = 1052.2222
$ws.Range("J110").Value = 4445
$ws.Range("K110").Value = 1052.2222
$ws.Range("L110").Value = 4445
$ws.Range("M110").Value = 992.7778000000001
$ws.Range("N110").Value = -8535
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = ""
$ws.Range("N113").Value = 0
$ws.Range("H116").Value = 882.1905
$ws.Range("I116").Value = 794.125
$ws.Range("J116").Value = 936.38464
$ws.Range("K116").Value = 794.125
$ws.Range("L116").Value = 936.38464
$ws.Range("M116").Value = 1499.875
$ws.Range("N116").Value = -5524.38464
$ws.Range("H122").Value = 2108.0454
$ws.Range("I122").Value = 2108.0454
$ws.Range("K122").Value = 6324.1362
$ws.Range("M122").Value = -3874.1362
$ws.Range("H132").Value = 4930.871
$ws.Range("I132").Value = 1384.36
$ws.Range("J132").Value = 19708
$ws.Range("K132").Value = 4153.08
$ws.Range("L132").Value = 59124
$ws.Range("M132").Value = -1623.08
$ws.Range("N132").Value = -64184

# --- BSM sheet updates ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 882.1905
$ws.Range("I3").Value = 794.125
$ws.Range("J3").Value = 936.38464
$ws.Range("K3").Value = 794.125
$ws.Range("L3").Value = 936.38464
$ws.Range("M3").Value = -680.125
$ws.Range("N3").Value = -1164.38464

# --- CRP sheet updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2126.5217
$ws.Range("I16").Value = 999.75
$ws.Range("J16").Value = 3355.7273
$ws.Range("K16").Value = 999.75
$ws.Range("L16").Value = 3355.7273
$ws.Range("M16").Value = -712.75
$ws.Range("N16").Value = -3929.7273
$ws.Range("H110").Value = 25000
$ws.Range("J110").Value = 25000
$ws.Range("L110").Value = 25000
$ws.Range("N110").Value = -33180
$ws.Range("H111").Value = 45702
$ws.Range("J111").Value = 45702
$ws.Range("L111").Value = 45702
$ws.Range("N111").Value = -53882
$ws.Range("H113").Value = 2126.5217
$ws.Range("I113").Value = 999.75
$ws.Range("J113").Value = 3355.7273
$ws.Range("K113").Value = 999.75
$ws.Range("L113").Value = 3355.7273
$ws.Range("M113").Value = 1170.25
$ws.Range("N113").Value = -7695.7273

# --- CUL sheet updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 8333842.5
$ws.Range("I113").Value = 5952845.5
$ws.Range("J113").Value = 12500587
$ws.Range("K113").Value = 17858536.5
$ws.Range("L113").Value = 37501761
$ws.Range("M113").Value = -17856366.5
$ws.Range("N113").Value = -37506101

# --- GSM sheet updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1171.2
$ws.Range("I102").Value = 1058.8572
$ws.Range("K102").Value = 1058.8572
$ws.Range("M102").Value = 563.1428000000001
$ws.Range("H122").Value = 7466097
$ws.Range("I122").Value = 10003471
$ws.Range("J122").Value = 3232.5293
$ws.Range("K122").Value = 30010413
$ws.Range("L122").Value = 9697.5879
$ws.Range("M122").Value = -30007963
$ws.Range("N122").Value = -14597.5879

# --- LTW sheet updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 35715384
$ws.Range("I40").Value = 1281.3334
$ws.Range("J40").Value = 250000000
$ws.Range("K40").Value = 1281.3334
$ws.Range("L40").Value = 250000000
$ws.Range("M40").Value = -1145.3334
$ws.Range("N40").Value = -250000272
$ws.Range("H132").Value = 6263.325
$ws.Range("I132").Value = 1870
$ws.Range("J132").Value = 10238.238
$ws.Range("K132").Value = 5610
$ws.Range("L132").Value = 30714.714
$ws.Range("M132").Value = -3080
$ws.Range("N132").Value = -35774.714

# --- WVR sheet updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 731.5714
$ws.Range("I100").Value = 813.8182
$ws.Range("J100").Value = 430
$ws.Range("K100").Value = 1627.6364
$ws.Range("L100").Value = 860
$ws.Range("M100").Value = -1086.6364
$ws.Range("N100").Value = -1942
$ws.Range("H108").Value = 36648
$ws.Range("J108").Value = 36648
$ws.Range("L108").Value = 36648
$ws.Range("N108").Value = -44328
$ws.Range("H122").Value = 1514.9412
$ws.Range("I122").Value = 1012.95654
$ws.Range("J122").Value = 2564.5454
$ws.Range("K122").Value = 3038.86962
$ws.Range("L122").Value = 7693.6362
$ws.Range("M122").Value = -588.8696199999999
$ws.Range("N122").Value = -12593.6362
